$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Price text has a numerically-insignificant trailing zero
# ("38.30", "3.50", "0.930") must be forced to Text format first, otherwise
# the normal cell-input parsing would coerce them to plain numbers
# (38.3, 3.5, 0.93) and the trailing zero would be lost.
foreach ($addr in @("D10", "D19", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '52.105.62'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '2.971.03'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '354.36'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").Value = '106.93'
$ws.Range("E6").Value = '  -4.72%  '
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.612'
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("D10").Value = '38.30'
$ws.Range("E10").Value = '  -2.91%  '
$ws.Range("D11").Value = '0.139'
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("E12").Value = '  -3.96%  '
$ws.Range("E13").Value = '  -4.16%  '
$ws.Range("D14").Value = '3.444.19'
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("D15").Value = '7.62'
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("D16").Value = '2.969.27'
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("E17").Value = '  +1.90%  '
$ws.Range("D18").Value = '52.123.18'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").Value = '3.50'
$ws.Range("E19").Value = '  +5.48%  '
$ws.Range("E20").Value = '  -2.25%  '
$ws.Range("D21").Value = '13.59'
$ws.Range("E21").Value = '  -4.47%  '
$ws.Range("D22").Value = '0.0₃0972'
$ws.Range("E22").Value = '  -1.42%  '
$ws.Range("D23").Value = '69.53'
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("D24").Value = '263.51'
$ws.Range("E24").Value = '  -1.99%  '
$ws.Range("E25").Value = '  -1.87%  '
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").Value = '26.79'
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").Value = '7.54'
$ws.Range("E28").Value = '  +2.45%  '
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").Value = '0.109'
$ws.Range("E30").Value = '  +1.78%  '
$ws.Range("E31").Value = '  -2.91%  '
$ws.Range("D32").Value = '6.09'
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("E33").Value = '  -3.61%  '
$ws.Range("E34").Value = '  -4.56%  '
$ws.Range("D35").Value = '50.76'
$ws.Range("E35").Value = '  -4.04%  '
$ws.Range("E36").Value = '  -2.14%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  -3.42%  '
$ws.Range("D39").Value = '17.88'
$ws.Range("E39").Value = '  -5.33%  '
$ws.Range("E40").Value = '  -4.32%  '
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("D43").Value = '22.64'
$ws.Range("E43").Value = '  -2.74%  '
$ws.Range("D44").Value = '122.01'
$ws.Range("E44").Value = '  +9.35%  '
$ws.Range("D45").Value = '2.12'
$ws.Range("E45").Value = '  -3.67%  '
$ws.Range("D46").Value = '2.115.47'
$ws.Range("E46").Value = '  -2.58%  '
$ws.Range("E47").Value = '  -4.24%  '
$ws.Range("D48").Value = '2.33'
$ws.Range("E48").Value = '  -8.04%  '
$ws.Range("D49").Value = '0.239'
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("D50").Value = '0.0335'
$ws.Range("E50").Value = '  -3.07%  '
$ws.Range("D51").Value = '0.930'
$ws.Range("E51").Value = '  -0.84%  '
